$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header date moved forward one month: 45406 (2024-04-24) -> 45436 (2024-05-24)
$ws.Range("A1").Value = "5/24/2024"

# Updated per-meter prices in the "SOGA de Monofilamento" price table
$ws.Range("D14").Value = 98.8
$ws.Range("D15").Value = 142
$ws.Range("D16").Value = 202
$ws.Range("D17").Value = 361
